$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.46659951392083
$ws.Range("D2").Value = 0.6453739777789236

$ws.Range("C3").Value = 0.9965854368087912
$ws.Range("D3").Value = 0.3298016152798811

$ws.Range("C4").Value = -0.1441525967470161
$ws.Range("D4").Value = 0.8866923300557712

$ws.Range("C5").Value = 0.2452526349988132
$ws.Range("D5").Value = 0.8085341228352818

$ws.Range("C6").Value = 0.3462654408956523
$ws.Range("D6").Value = 0.7324329763674762

$ws.Range("C7").Value = -0.8770656631997783
$ws.Range("D7").Value = 0.3899282544463496

$ws.Range("C8").Value = -0.2870552726176914
$ws.Range("D8").Value = 0.7767562845370648

$ws.Range("C9").Value = -1.152750111812624
$ws.Range("D9").Value = 0.2613849018011334

$ws.Range("C10").Value = -0.5564949194574234
$ws.Range("D10").Value = 0.5834894822473535

$ws.Range("C11").Value = 0.4173619053608623
$ws.Range("D11").Value = 0.6804585857483167
